# Edit workbook cotisations_TousRegimes.xlsx:
# - add two new aggregate rows "cotsoc_noncontrib" / "cotsoc_contrib" to
#   sheet "amounts" (sheet1) and sheet "montants" (sheet2, with formulas).

$wb  = $excel.ActiveWorkbook
$wsA = $wb.Worksheets.Item("amounts")
$wsM = $wb.Worksheets.Item("montants")

# ---------------------------------------------------------------------
# 1) "montants" sheet: insert two new rows right after existing row 8
#    (the "cotisation des non-salaries" row), pushing everything below
#    down by two rows, then fill them in with the new aggregate data.
# ---------------------------------------------------------------------
$wsM.Rows("9:10").Insert()

# Copy the formatting of the row immediately above (row 8) onto the two
# freshly inserted rows so they pick up the same number formats/borders.
$wsM.Range("A8:I8").Copy()
$wsM.Range("A9:I10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$wsM.Cells.Item(9, 1).Value = "cotsoc_noncontrib"
$wsM.Range("B9").Formula = "=B15+B17"
$wsM.Range("C9:I9").Formula = "=C15+C17"

$wsM.Cells.Item(10, 1).Value = "cotsoc_contrib"
$wsM.Range("B10").Formula = "=B16+B18"
$wsM.Range("C10:I10").Formula = "=C16+C18"

# ---------------------------------------------------------------------
# 2) "amounts" sheet: append the same two aggregate rows (plain values,
#    no formulas) right after the existing last row (row 7).
# ---------------------------------------------------------------------
# Row 5 already carries the "label in col A (border style) + numeric
# value in B:I (number-format style)" combination we need for the two
# new rows, so reuse its formatting.
$wsA.Range("A5:I5").Copy()
$wsA.Range("A8:I9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$wsA.Cells.Item(8, 1).Value = "cotsoc_noncontrib"
$noncontribValues = @(92806000000, 95002000000, 97978000000, 101486000000, 104001000000, 108053000000, 108818000000, 110926000000)
for ($i = 0; $i -lt 8; $i++) {
    $wsA.Cells.Item(8, 2 + $i).Value = $noncontribValues[$i]
}

$wsA.Cells.Item(9, 1).Value = "cotsoc_contrib"
$contribValues = @(89157000000, 92100000000, 94188000000, 98386000000, 102167000000, 104647000000, 105452000000, 112923000000)
for ($i = 0; $i -lt 8; $i++) {
    $wsA.Cells.Item(9, 2 + $i).Value = $contribValues[$i]
}

# ---------------------------------------------------------------------
# 3) View/selection tweaks described in the diff. "amounts" remains the
#    active/selected tab, so re-activate it last.
# ---------------------------------------------------------------------
$wsM.Activate()
$excel.ActiveWindow.ScrollRow = 7
$wsM.Range("A9:I10").Select()

$wsA.Activate()
$wsA.Range("D16").Select()

Write-Host "done"
